$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Typography")

$ws.Range("G5").Value = "°.%/"
$ws.Range("G6").Value = "°.%/"
$ws.Range("G7").Value = "°.%/"

$ws.Range("I5").Value = "a-z,A-Z,0-9"
$ws.Range("I6").Value = "a-z,A-Z,0-9"
$ws.Range("I7").Value = "a-z,A-Z,0-9"
